$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Enter values in the same sequence the original author used, so that the
# shared-strings table is rebuilt with a matching order of new entries.

# New TestCase Name labels for the added rows
$ws.Range("A4").Value = "TC_03"
$ws.Range("A5").Value = "TC_04"
$ws.Range("A6").Value = "TC_05"

# Amount column: existing rows updated (kept as quote-prefixed text, matching
# their original formatting), then new rows filled in
$ws.Range("B3").Value = "'90000"
$ws.Range("B2").Value = "'100"
$ws.Range("B4").Value = "85000"
$ws.Range("B5").Value = "100000"

# Source Currency for TC_04
$ws.Range("C5").Value = "USD"

# Remaining Amount entry
$ws.Range("B6").Value = "8797"

# Remaining currency cells (reuse already-existing shared strings)
$ws.Range("D2").Value = "GBP"
$ws.Range("D3").Value = "INR"
$ws.Range("C4").Value = "GBP"
$ws.Range("D4").Value = "INR"
$ws.Range("D5").Value = "EUR"
$ws.Range("C6").Value = "EUR"
$ws.Range("D6").Value = "USD"

$ws.Range("D6").Select()
